# Verifiche documenti.xlsx
# Commit message: "sistemata larghezza colonne che era saltata, tolto commento
# tra parentesi" (fixed column widths that had jumped/been lost, removed the
# parenthetical comment)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# ------------------------------------------------------------------
# 1) "tolto commento tra parentesi": E10 used to hold the text
#    "39(contando parole senza g)" - replace it with the plain number 39
#    (this also drops that string from the shared-strings table).
# ------------------------------------------------------------------
$ws.Range("E10").Value = 39

# B11 (value 126.30) picks up the same right-aligned numeric style used by
# E10/E11.
$ws.Range("B11").HorizontalAlignment = -4152   # xlRight

# Header row (B1:G1, bold style) gets centered alignment added.
$ws.Range("B1:G1").HorizontalAlignment = -4108  # xlCenter

# ------------------------------------------------------------------
# 2) "sistemata larghezza colonne che era saltata": restore the column
#    widths that had been lost/jumped (target character widths: A=32.83,
#    B=C=9.66, D=E=8.83, F=9.16, G=9.5).
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 32
$ws.Columns.Item(2).ColumnWidth = 53 / 6
$ws.Columns.Item(3).ColumnWidth = 53 / 6
$ws.Columns.Item(4).ColumnWidth = 8
$ws.Columns.Item(5).ColumnWidth = 8
$ws.Columns.Item(6).ColumnWidth = 50 / 6
$ws.Columns.Item(7).ColumnWidth = 52 / 6

# Row 1 height adjusts slightly along with the other view refresh.
$ws.Rows.Item(1).RowHeight = 19

# Selection moves from D11 to B1.
$null = $ws.Range("B1").Select()
